$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.988.65"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "'1.819.94"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'309.69"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "'0.4633"
$ws.Range("E7").Value = "  -2.63%  "
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "'0.07280"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "'0.8653"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "'1.884.49"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'0.07606"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").Value = "'93.07"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'5.329"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "'6.468"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'0.000008631"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "'27.409.12"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'14.47"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").Value = "'5.156"
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").Value = "'10.58"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").Value = "'2.108.64"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "'151.73"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").Value = "'1.859"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").Value = "'18.22"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "'2.097"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").Value = "'5.082"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "'115.88"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "'0.08898"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'2.952"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "'0.7314"
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").Value = "'4.429"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("D36").Value = "'1.009"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'2.533"
$ws.Range("E37").Value = "  +6.36%  "
$ws.Range("D38").Value = "'1.076"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").Value = "'0.05264"
$ws.Range("D40").Value = "'0.01916"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").Value = "'2.935"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("D42").Value = "'7.136"
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("D43").Value = "'0.5211"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "'0.1633"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'8.237"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("D46").Value = "'0.4851"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "'10.10"
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("D49").Value = "'103.16"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("D50").Value = "'1.634"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").Value = "'0.06224"
$ws.Range("E51").Value = "  -1.59%  "
